# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamp cells for the
# 85292b46-8ed0-4aa3-815a-34da459a2008.md row across the Overview, zh-cn and
# de-de sheets, reflecting a freshly (re)generated handback report.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the 85292b46... file
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-21 22:54:45"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 85292b46... file
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-21 22:54:40"
$wsZhCn.Range("K4").Value = "2016-08-21 22:54:57"

# de-de sheet: "Correspond Handback DateTime" for the 85292b46... file
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-21 22:55:08"
